$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2 to hold the string value "365652-7" instead of the numeric 365652
$ws.Range("D2").Value = "365652-7"

# Update the active cell selection to D3
$ws.Range("D3").Select()
